{"js": "const title = { oldText: \"2024-03-05 Tuesday\", newText: \"2024-03-06 Wednesday\" };\nconst cellEdits = [\n  { row: 0, col: 0, oldText: \"60-16=\", newText: \"91-48=\" },\n  { row: 0, col: 1, oldText: \"48+28=\", newText: \"38+27=\" },\n  { row: 0, col: 2, oldText: \"6+68=\", newText: \"63-36=\" },\n  { row: 0, col: 3, oldText: \"80-44=\", newText: \"60-26=\" },\n  { row: 0, col: 4, oldText: \"4+87=\", newText: \"51-27=\" },\n  { row: 1, col: 0, oldText: \"68+14=\", newText: \"57+28=\" },\n  { row: 1, col: 1, oldText: \"67+9=\", newText: \"73+9=\" },\n  { row: 1, col: 2, oldText: \"96-77=\", newText: \"74-48=\" },\n  { row: 1, col: 3, oldText: \"17-8=\", newText: \"53-18=\" },\n  { row: 1, col: 4, oldText: \"42+49=\", newText: \"51-32=\" },\n  { row: 2, col: 0, oldText: \"51-22=\", newText: \"70-4=\" },\n  { row: 2, col: 1, oldText: \"81-24=\", newText: \"8+29=\" },\n  { row: 2, col: 2, oldText: \"57-28=\", newText: \"93-88=\" },\n  { row: 2, col: 3, oldText: \"34-17=\", newText: \"5+7=\" },\n  { row: 2, col: 4, oldText: \"85-26=\", newText: \"73-4=\" },\n  { row: 3, col: 0, oldText: \"11-8=\", newText: \"10-6=\" },\n  { row: 3, col: 1, oldText: \"61-43=\", newText: \"63-5=\" },\n  { row: 3, col: 2, oldText: \"82-53=\", newText: \"15-6=\" },\n  { row: 3, col: 3, oldText: \"70-46=\", newText: \"81-79=\" },\n  { row: 3, col: 4, oldText: \"56+39=\", newText: \"67+8=\" },\n  { row: 4, col: 0, oldText: \"38+45=\", newText: \"70-34=\" },\n  { row: 4, col: 1, oldText: \"64-27=\", newText: \"64+19=\" },\n  { row: 4, col: 2, oldText: \"15+76=\", newText: \"55+16=\" },\n  { row: 4, col: 3, oldText: \"75-57=\", newText: \"36+28=\" },\n  { row: 4, col: 4, oldText: \"14+78=\", newText: \"55-36=\" },\n  { row: 5, col: 0, oldText: \"93-76=\", newText: \"32+9=\" },\n  { row: 5, col: 1, oldText: \"13-9=\", newText: \"80-18=\" },\n  { row: 5, col: 2, oldText: \"8+76=\", newText: \"58+39=\" },\n  { row: 5, col: 3, oldText: \"74+8=\", newText: \"37-9=\" },\n  { row: 5, col: 4, oldText: \"87+8=\", newText: \"67+24=\" },\n  { row: 6, col: 0, oldText: \"73-39=\", newText: \"70-36=\" },\n  { row: 6, col: 1, oldText: \"41-24=\", newText: \"80-48=\" },\n  { row: 6, col: 2, oldText: \"84-65=\", newText: \"94-45=\" },\n  { row: 6, col: 3, oldText: \"68-49=\", newText: \"32-28=\" },\n  { row: 6, col: 4, oldText: \"70-7=\", newText: \"30-14=\" },\n  { row: 7, col: 0, oldText: \"64-55=\", newText: \"27+16=\" },\n  { row: 7, col: 1, oldText: \"50-9=\", newText: \"44+39=\" },\n  { row: 7, col: 2, oldText: \"70-43=\", newText: \"71-37=\" },\n  { row: 7, col: 3, oldText: \"23-4=\", newText: \"77+5=\" },\n  { row: 7, col: 4, oldText: \"73-28=\", newText: \"50-27=\" },\n  { row: 8, col: 0, oldText: \"69+7=\", newText: \"15-9=\" },\n  { row: 8, col: 1, oldText: \"96-29=\", newText: \"66+25=\" },\n  { row: 8, col: 2, oldText: \"68+9=\", newText: \"27-8=\" },\n  { row: 8, col: 3, oldText: \"24+58=\", newText: \"96-89=\" },\n  { row: 8, col: 4, oldText: \"19+58=\", newText: \"9+33=\" },\n  { row: 9, col: 0, oldText: \"93-79=\", newText: \"46+16=\" },\n  { row: 9, col: 1, oldText: \"69+23=\", newText: \"92-86=\" },\n  { row: 9, col: 2, oldText: \"25+46=\", newText: \"19+38=\" },\n  { row: 9, col: 3, oldText: \"12-9=\", newText: \"29+54=\" },\n  { row: 9, col: 4, oldText: \"39+48=\", newText: \"76-28=\" },\n  { row: 10, col: 0, oldText: \"98-9=\", newText: \"48+9=\" },\n  { row: 10, col: 1, oldText: \"64+8=\", newText: \"82-44=\" },\n  { row: 10, col: 2, oldText: \"35+28=\", newText: \"37+24=\" },\n  { row: 10, col: 3, oldText: \"7+58=\", newText: \"47+26=\" },\n  { row: 10, col: 4, oldText: \"70-36=\", newText: \"36+19=\" },\n  { row: 11, col: 0, oldText: \"84-16=\", newText: \"53+18=\" },\n  { row: 11, col: 1, oldText: \"94-35=\", newText: \"18+56=\" },\n  { row: 11, col: 2, oldText: \"36+6=\", newText: \"49+22=\" },\n  { row: 11, col: 3, oldText: \"51-8=\", newText: \"90-69=\" },\n  { row: 11, col: 4, oldText: \"12+39=\", newText: \"4+79=\" },\n  { row: 12, col: 0, oldText: \"30-24=\", newText: \"71-33=\" },\n  { row: 12, col: 1, oldText: \"93-85=\", newText: \"25+16=\" },\n  { row: 12, col: 2, oldText: \"84-15=\", newText: \"66-38=\" },\n  { row: 12, col: 3, oldText: \"25+57=\", newText: \"72+19=\" },\n  { row: 12, col: 4, oldText: \"42+29=\", newText: \"93-6=\" },\n  { row: 13, col: 0, oldText: \"23+59=\", newText: \"63-58=\" },\n  { row: 13, col: 1, oldText: \"75-68=\", newText: \"27+54=\" },\n  { row: 13, col: 2, oldText: \"91-64=\", newText: \"81-56=\" },\n  { row: 13, col: 3, oldText: \"5+9=\", newText: \"60-47=\" },\n  { row: 13, col: 4, oldText: \"76+18=\", newText: \"63+9=\" },\n  { row: 14, col: 0, oldText: \"28+65=\", newText: \"22+19=\" },\n  { row: 14, col: 1, oldText: \"32+19=\", newText: \"24+29=\" },\n  { row: 14, col: 2, oldText: \"51-23=\", newText: \"26+36=\" },\n  { row: 14, col: 3, oldText: \"9+49=\", newText: \"60-1=\" },\n  { row: 14, col: 4, oldText: \"34+38=\", newText: \"95-48=\" },\n  { row: 15, col: 0, oldText: \"38+4=\", newText: \"14+19=\" },\n  { row: 15, col: 1, oldText: \"28+25=\", newText: \"38+47=\" },\n  { row: 15, col: 2, oldText: \"76-8=\", newText: \"21-12=\" },\n  { row: 15, col: 3, oldText: \"67-59=\", newText: \"8+38=\" },\n  { row: 15, col: 4, oldText: \"27+25=\", newText: \"46+29=\" },\n  { row: 16, col: 0, oldText: \"92-17=\", newText: \"70-45=\" },\n  { row: 16, col: 1, oldText: \"31-2=\", newText: \"71-17=\" },\n  { row: 16, col: 2, oldText: \"9+42=\", newText: \"15+66=\" },\n  { row: 16, col: 3, oldText: \"67-19=\", newText: \"70-28=\" },\n  { row: 16, col: 4, oldText: \"48+17=\", newText: \"57-29=\" },\n  { row: 17, col: 0, oldText: \"62-13=\", newText: \"17+7=\" },\n  { row: 17, col: 1, oldText: \"16+8=\", newText: \"65+28=\" },\n  { row: 17, col: 2, oldText: \"93-17=\", newText: \"11-4=\" },\n  { row: 17, col: 3, oldText: \"61-26=\", newText: \"19+52=\" },\n  { row: 17, col: 4, oldText: \"18+8=\", newText: \"14+9=\" },\n  { row: 18, col: 0, oldText: \"43-8=\", newText: \"69+2=\" },\n  { row: 18, col: 1, oldText: \"95-18=\", newText: \"35-27=\" },\n  { row: 18, col: 2, oldText: \"31-7=\", newText: \"54-29=\" },\n  { row: 18, col: 3, oldText: \"26+26=\", newText: \"76-48=\" },\n  { row: 18, col: 4, oldText: \"70-27=\", newText: \"73-34=\" },\n  { row: 19, col: 0, oldText: \"70-21=\", newText: \"83-5=\" },\n  { row: 19, col: 1, oldText: \"19+62=\", newText: \"48+25=\" },\n  { row: 19, col: 2, oldText: \"28+65=\", newText: \"17+48=\" },\n  { row: 19, col: 3, oldText: \"29+34=\", newText: \"36+59=\" },\n  { row: 19, col: 4, oldText: \"83-35=\", newText: \"70-35=\" },\n];\n// --- Update the title line (first paragraph of the body, above the table). ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items.find((p) => p.text === title.oldText);\nif (!titlePara) {\n  throw new Error(\"Could not find title paragraph with text: \" + title.oldText);\n}\ntitlePara.getRange().insertText(title.newText, \"Replace\");\n\n// --- Update every math-problem cell in the table by (row, col) position, ---\n// verifying the existing text before replacing it (defensive; a couple of\n// source strings like \"28+65=\" repeat at different grid positions with\n// different replacements, so matching purely by text is unsafe). Cell\n// lookups/loads are batched into two sync round-trips instead of one per\n// cell.\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"isNullObject\");\nawait context.sync();\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\nconst cells = cellEdits.map((edit) => {\n  const cell = table.getCell(edit.row, edit.col);\n  cell.load(\"value\");\n  return cell;\n});\nawait context.sync();\n\ncellEdits.forEach((edit, i) => {\n  const cell = cells[i];\n  const current = cell.value.replace(/\\r|\\u0007/g, \"\");\n  if (current !== edit.oldText) {\n    throw new Error(\n      \"Cell (\" + edit.row + \",\" + edit.col + \") expected '\" + edit.oldText +\n      \"' but found '\" + current + \"'\"\n    );\n  }\n  cell.value = edit.newText;\n});\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Update the title line (first paragraph, above the table). ---\n$titleOld = '2024-03-05 Tuesday'\n$titleNew = '2024-03-06 Wednesday'\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($titleText -ne $titleOld) {\n    throw \"Title paragraph expected $titleOld but found $titleText\"\n}\n$titlePara.Range.Text = $titleNew\n\n# --- Update every math-problem cell in the table by (row, col), ---\n# --- verifying the existing text before replacing it (some source ---\n# --- strings repeat at different grid positions with different ---\n# --- replacements, so a blind global Find/Replace is not safe). ---\n$table = $d.Tables.Item(1)\n$edits = @(\n    @{ Row = 1; Col = 1; OldText = '60-16='; NewText = '91-48=' },\n    @{ Row = 1; Col = 2; OldText = '48+28='; NewText = '38+27=' },\n    @{ Row = 1; Col = 3; OldText = '6+68='; NewText = '63-36=' },\n    @{ Row = 1; Col = 4; OldText = '80-44='; NewText = '60-26=' },\n    @{ Row = 1; Col = 5; OldText = '4+87='; NewText = '51-27=' },\n    @{ Row = 2; Col = 1; OldText = '68+14='; NewText = '57+28=' },\n    @{ Row = 2; Col = 2; OldText = '67+9='; NewText = '73+9=' },\n    @{ Row = 2; Col = 3; OldText = '96-77='; NewText = '74-48=' },\n    @{ Row = 2; Col = 4; OldText = '17-8='; NewText = '53-18=' },\n    @{ Row = 2; Col = 5; OldText = '42+49='; NewText = '51-32=' },\n    @{ Row = 3; Col = 1; OldText = '51-22='; NewText = '70-4=' },\n    @{ Row = 3; Col = 2; OldText = '81-24='; NewText = '8+29=' },\n    @{ Row = 3; Col = 3; OldText = '57-28='; NewText = '93-88=' },\n    @{ Row = 3; Col = 4; OldText = '34-17='; NewText = '5+7=' },\n    @{ Row = 3; Col = 5; OldText = '85-26='; NewText = '73-4=' },\n    @{ Row = 4; Col = 1; OldText = '11-8='; NewText = '10-6=' },\n    @{ Row = 4; Col = 2; OldText = '61-43='; NewText = '63-5=' },\n    @{ Row = 4; Col = 3; OldText = '82-53='; NewText = '15-6=' },\n    @{ Row = 4; Col = 4; OldText = '70-46='; NewText = '81-79=' },\n    @{ Row = 4; Col = 5; OldText = '56+39='; NewText = '67+8=' },\n    @{ Row = 5; Col = 1; OldText = '38+45='; NewText = '70-34=' },\n    @{ Row = 5; Col = 2; OldText = '64-27='; NewText = '64+19=' },\n    @{ Row = 5; Col = 3; OldText = '15+76='; NewText = '55+16=' },\n    @{ Row = 5; Col = 4; OldText = '75-57='; NewText = '36+28=' },\n    @{ Row = 5; Col = 5; OldText = '14+78='; NewText = '55-36=' },\n    @{ Row = 6; Col = 1; OldText = '93-76='; NewText = '32+9=' },\n    @{ Row = 6; Col = 2; OldText = '13-9='; NewText = '80-18=' },\n    @{ Row = 6; Col = 3; OldText = '8+76='; NewText = '58+39=' },\n    @{ Row = 6; Col = 4; OldText = '74+8='; NewText = '37-9=' },\n    @{ Row = 6; Col = 5; OldText = '87+8='; NewText = '67+24=' },\n    @{ Row = 7; Col = 1; OldText = '73-39='; NewText = '70-36=' },\n    @{ Row = 7; Col = 2; OldText = '41-24='; NewText = '80-48=' },\n    @{ Row = 7; Col = 3; OldText = '84-65='; NewText = '94-45=' },\n    @{ Row = 7; Col = 4; OldText = '68-49='; NewText = '32-28=' },\n    @{ Row = 7; Col = 5; OldText = '70-7='; NewText = '30-14=' },\n    @{ Row = 8; Col = 1; OldText = '64-55='; NewText = '27+16=' },\n    @{ Row = 8; Col = 2; OldText = '50-9='; NewText = '44+39=' },\n    @{ Row = 8; Col = 3; OldText = '70-43='; NewText = '71-37=' },\n    @{ Row = 8; Col = 4; OldText = '23-4='; NewText = '77+5=' },\n    @{ Row = 8; Col = 5; OldText = '73-28='; NewText = '50-27=' },\n    @{ Row = 9; Col = 1; OldText = '69+7='; NewText = '15-9=' },\n    @{ Row = 9; Col = 2; OldText = '96-29='; NewText = '66+25=' },\n    @{ Row = 9; Col = 3; OldText = '68+9='; NewText = '27-8=' },\n    @{ Row = 9; Col = 4; OldText = '24+58='; NewText = '96-89=' },\n    @{ Row = 9; Col = 5; OldText = '19+58='; NewText = '9+33=' },\n    @{ Row = 10; Col = 1; OldText = '93-79='; NewText = '46+16=' },\n    @{ Row = 10; Col = 2; OldText = '69+23='; NewText = '92-86=' },\n    @{ Row = 10; Col = 3; OldText = '25+46='; NewText = '19+38=' },\n    @{ Row = 10; Col = 4; OldText = '12-9='; NewText = '29+54=' },\n    @{ Row = 10; Col = 5; OldText = '39+48='; NewText = '76-28=' },\n    @{ Row = 11; Col = 1; OldText = '98-9='; NewText = '48+9=' },\n    @{ Row = 11; Col = 2; OldText = '64+8='; NewText = '82-44=' },\n    @{ Row = 11; Col = 3; OldText = '35+28='; NewText = '37+24=' },\n    @{ Row = 11; Col = 4; OldText = '7+58='; NewText = '47+26=' },\n    @{ Row = 11; Col = 5; OldText = '70-36='; NewText = '36+19=' },\n    @{ Row = 12; Col = 1; OldText = '84-16='; NewText = '53+18=' },\n    @{ Row = 12; Col = 2; OldText = '94-35='; NewText = '18+56=' },\n    @{ Row = 12; Col = 3; OldText = '36+6='; NewText = '49+22=' },\n    @{ Row = 12; Col = 4; OldText = '51-8='; NewText = '90-69=' },\n    @{ Row = 12; Col = 5; OldText = '12+39='; NewText = '4+79=' },\n    @{ Row = 13; Col = 1; OldText = '30-24='; NewText = '71-33=' },\n    @{ Row = 13; Col = 2; OldText = '93-85='; NewText = '25+16=' },\n    @{ Row = 13; Col = 3; OldText = '84-15='; NewText = '66-38=' },\n    @{ Row = 13; Col = 4; OldText = '25+57='; NewText = '72+19=' },\n    @{ Row = 13; Col = 5; OldText = '42+29='; NewText = '93-6=' },\n    @{ Row = 14; Col = 1; OldText = '23+59='; NewText = '63-58=' },\n    @{ Row = 14; Col = 2; OldText = '75-68='; NewText = '27+54=' },\n    @{ Row = 14; Col = 3; OldText = '91-64='; NewText = '81-56=' },\n    @{ Row = 14; Col = 4; OldText = '5+9='; NewText = '60-47=' },\n    @{ Row = 14; Col = 5; OldText = '76+18='; NewText = '63+9=' },\n    @{ Row = 15; Col = 1; OldText = '28+65='; NewText = '22+19=' },\n    @{ Row = 15; Col = 2; OldText = '32+19='; NewText = '24+29=' },\n    @{ Row = 15; Col = 3; OldText = '51-23='; NewText = '26+36=' },\n    @{ Row = 15; Col = 4; OldText = '9+49='; NewText = '60-1=' },\n    @{ Row = 15; Col = 5; OldText = '34+38='; NewText = '95-48=' },\n    @{ Row = 16; Col = 1; OldText = '38+4='; NewText = '14+19=' },\n    @{ Row = 16; Col = 2; OldText = '28+25='; NewText = '38+47=' },\n    @{ Row = 16; Col = 3; OldText = '76-8='; NewText = '21-12=' },\n    @{ Row = 16; Col = 4; OldText = '67-59='; NewText = '8+38=' },\n    @{ Row = 16; Col = 5; OldText = '27+25='; NewText = '46+29=' },\n    @{ Row = 17; Col = 1; OldText = '92-17='; NewText = '70-45=' },\n    @{ Row = 17; Col = 2; OldText = '31-2='; NewText = '71-17=' },\n    @{ Row = 17; Col = 3; OldText = '9+42='; NewText = '15+66=' },\n    @{ Row = 17; Col = 4; OldText = '67-19='; NewText = '70-28=' },\n    @{ Row = 17; Col = 5; OldText = '48+17='; NewText = '57-29=' },\n    @{ Row = 18; Col = 1; OldText = '62-13='; NewText = '17+7=' },\n    @{ Row = 18; Col = 2; OldText = '16+8='; NewText = '65+28=' },\n    @{ Row = 18; Col = 3; OldText = '93-17='; NewText = '11-4=' },\n    @{ Row = 18; Col = 4; OldText = '61-26='; NewText = '19+52=' },\n    @{ Row = 18; Col = 5; OldText = '18+8='; NewText = '14+9=' },\n    @{ Row = 19; Col = 1; OldText = '43-8='; NewText = '69+2=' },\n    @{ Row = 19; Col = 2; OldText = '95-18='; NewText = '35-27=' },\n    @{ Row = 19; Col = 3; OldText = '31-7='; NewText = '54-29=' },\n    @{ Row = 19; Col = 4; OldText = '26+26='; NewText = '76-48=' },\n    @{ Row = 19; Col = 5; OldText = '70-27='; NewText = '73-34=' },\n    @{ Row = 20; Col = 1; OldText = '70-21='; NewText = '83-5=' },\n    @{ Row = 20; Col = 2; OldText = '19+62='; NewText = '48+25=' },\n    @{ Row = 20; Col = 3; OldText = '28+65='; NewText = '17+48=' },\n    @{ Row = 20; Col = 4; OldText = '29+34='; NewText = '36+59=' },\n    @{ Row = 20; Col = 5; OldText = '83-35='; NewText = '70-35=' }\n)\n\nforeach ($edit in $edits) {\n    $cell = $table.Cell($edit.Row, $edit.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $edit.OldText) {\n        throw \"Cell ($($edit.Row),$($edit.Col)) expected $($edit.OldText) but found $current\"\n    }\n    $cell.Range.Text = $edit.NewText\n}"}
